$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Cod cliente 5802202) so the row below (19499545) shifts up
$ws.Rows.Item(2).Delete()

# Update the active selection to C8, matching the post-edit workbook state
$ws.Range("C8").Select()
